# Refresh the cryptos price/volume list (GitHub Actions scheduled update).
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's existing
# inline-string cells) instead of re-parsing/reformatting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.708.59"
$ws.Range("E2").Value = "  -4.91%  "

$ws.Range("D3").Value = "3.166.88"
$ws.Range("E3").Value = "  -5.23%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'531.65"
$ws.Range("E5").Value = "  -6.12%  "

$ws.Range("D6").Value = "'134.69"
$ws.Range("E6").Value = "  -7.95%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.169.50"
$ws.Range("E8").Value = "  -5.09%  "

$ws.Range("D9").Value = "'0.452"
$ws.Range("E9").Value = "  -6.40%  "

$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "  -8.05%  "

$ws.Range("E11").Value = "  -8.53%  "

$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "  -4.90%  "

$ws.Range("D13").Value = "3.715.38"

$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "'25.86"
$ws.Range("E15").Value = "  -6.48%  "

$ws.Range("D16").Value = "3.174.81"
$ws.Range("E16").Value = "  -5.24%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000154"
$ws.Range("E17").Value = "  -8.41%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "57.881.01"
$ws.Range("E18").Value = "  -4.64%  "

$ws.Range("D19").Value = "'5.84"
$ws.Range("E19").Value = "  -6.80%  "

$ws.Range("E20").Value = "  -8.63%  "

$ws.Range("E21").Value = "  -9.01%  "

$ws.Range("D22").Value = "'357.86"
$ws.Range("E22").Value = "  -4.79%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Value = "'69.48"
$ws.Range("E24").Value = "  -7.04%  "

$ws.Range("D25").Value = "'0.517"
$ws.Range("E25").Value = "  -7.62%  "

$ws.Range("D26").Value = "3.309.88"
$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("E27").Value = "  -3.42%  "

$ws.Range("D28").Value = "0.0₃0958"
$ws.Range("E28").Value = "  -11.00%  "

$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").Value = "'6.94"

$ws.Range("E32").Value = "  -8.82%  "

$ws.Range("E33").Value = "  -9.69%  "

$ws.Range("D34").Value = "'21.71"
$ws.Range("E34").Value = "  -4.88%  "

$ws.Range("E35").Value = "  -6.53%  "

$ws.Range("D36").Value = "'4.97"
$ws.Range("E36").Value = "  -6.52%  "

$ws.Range("D37").Value = "'160.24"
$ws.Range("E37").Value = "  -4.84%  "

$ws.Range("E38").Value = "  -7.50%  "

$ws.Range("D39").Value = "'6.28"
$ws.Range("E39").Value = "  -7.86%  "

$ws.Range("D40").Value = "'25.86"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("D41").Value = "'0.0705"
$ws.Range("E41").Value = "  -5.91%  "

$ws.Range("D42").Value = "3.200.30"

$ws.Range("E43").Value = "  -4.64%  "

$ws.Range("D44").Value = "'0.703"
$ws.Range("E44").Value = "  -7.09%  "

$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("E46").Value = "  -6.87%  "

$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  -7.84%  "

$ws.Range("D49").Value = "2.273.31"
$ws.Range("E49").Value = "  -7.52%  "

$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("D51").Value = "'20.71"
$ws.Range("E51").Value = "  -6.92%  "
